$wb = $excel.ActiveWorkbook

# --- "Air" sheet: translate cause labels from Indonesian to English ---
$air = $wb.Worksheets.Item("Air")
$air.Range("C2").Value = "Transportation"
$air.Range("C3").Value = "Industry"
$air.Range("C4").Value = "Domestic"
$air.Range("C5").Value = "PP & Heater"
$air.Range("A5").Value = "Power Plant & Heater"
$air.Range("A4").Value = "Domestic Burning"
$air.Range("A3").Value = "Industrial Smoke"

# --- "Car" sheet: translate fuel type ("Jenis") labels from Indonesian to English ---
$car = $wb.Worksheets.Item("Car")
$car.Range("E7").Value = "Elektric"
$car.Range("E8").Value = "Elektric"
$car.Range("E9").Value = "Elektric"
$car.Range("E10").Value = "Elektric"
$car.Range("E11").Value = "Elektric"
$car.Range("E2").Value = "Gasoline"
$car.Range("E3").Value = "Gasoline"
$car.Range("E4").Value = "Gasoline"
$car.Range("E5").Value = "Gasoline"
$car.Range("E6").Value = "Gasoline"
$car.Range("E12").Value = "Gasoline"

# --- Restore per-sheet cursor/selection positions ---
$norway = $wb.Worksheets.Item("Norway")
$norway.Range("C12").Select()

$air.Range("C15").Select()

# --- "Car" becomes the active/selected tab, cursor at C13 ---
$car.Activate()
$car.Range("C13").Select()
